$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the "title" column (B) so it can fit the Google Drive image/links
# that are now being dropped into it (~53.57 chars wide in the saved file).
$ws.Columns("B").ColumnWidth = 52.7

# The cursor ended up on O5 when the file was saved.
$ws.Range("O5").Select()

# Touch the fill of the last header cell (F1) so its cached style entry gets
# rebuilt -- on the real file this cell's xf lost its redundant
# applyFill="1" flag (fillId was already 0/"no fill"). Flipping the pattern
# away and back forces the engine to recompute that cell's style record.
$f1 = $ws.Range("F1")
$f1.Interior.Pattern = 17
$f1.Interior.Pattern = -4142
